{"js": "// Remove the word \"la \" from the phrase so it reads:\n// \"...solucionar problemas relacionados con la inform\u00e1tica\"\n// -> \"...solucionar problemas relacionados con inform\u00e1tica\"\nconst searchResults = context.document.body.search(\"relacionados con la \", {\n  matchCase: true,\n  matchWholeWord: false\n});\nsearchResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].insertText(\"relacionados con \", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# \"Se me da bien generar ideas que permitan solucionar problemas relacionados\n# con la inform\u00e1tica\" -> \"...relacionados con inform\u00e1tica\"\n# (drop the word \"la \" right before \"inform\u00e1tica\")\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"relacionados con la \"\n$find.Replacement.Text = \"relacionados con \"\n$find.Forward = $true\n$find.Wrap = 0\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 1)\n"}
